$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New table rows (36-42) -------------------------------------------------
# Entered in the same order the shared-string table shows they were first
# authored in (check-rows first, then the enter-rows slotted back in between,
# then the two price-option rows).

# Row 36: check hints (mandatory fields) - Product Page
$ws.Range("A36").Value = "Product Page check for hints regarding mandatory fields"
$ws.Range("B36").Value = "<CHK>"
$ws.Range("E36").Value = "Product Page check for hints regarding mandatory fields"
$ws.Range("H36").Value = "<NOP>"

# Row 38: check hint date with invalid format
$ws.Range("A38").Value = "Product Page check for hint date with invalid format"
$ws.Range("B38").Value = "<CHK>"
$ws.Range("E38").Value = "Product Page check for hint date with invalid format"
$ws.Range("H38").Value = "<NOP>"

# Row 40: check hint date with invalid value in past
$ws.Range("A40").Value = "Product Page check for hint date with invalid value in past"
$ws.Range("B40").Value = "<CHK>"
$ws.Range("E40").Value = "Product Page check for hint date with invalid value in past"
$ws.Range("H40").Value = "<NOP>"

# Row 37: enter date with invalid format
$ws.Range("A37").Value = "Product Page enter date with invalid format"
$ws.Range("A37").NumberFormat = "@"
$ws.Range("B37").Value = "<SET>"
$ws.Range("E37").Value = "Product Page enter date with invalid format"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("H37").Value = "<NOP>"

# Row 39: enter date with invalid value in past
$ws.Range("A39").Value = "Product Page enter date with invalid value in past"
$ws.Range("A39").NumberFormat = "@"
$ws.Range("B39").Value = "<SET>"
$ws.Range("E39").Value = "Product Page enter date with invalid value in past"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("H39").Value = "<NOP>"

# Row 41: price option page - open mandatory field
$ws.Range("A41").Value = "Price option page check for open mandatory field"
$ws.Range("B41").Value = "<CHK>"
$ws.Range("F41").Value = "Price option page check for open mandatory field"
$ws.Range("H41").Value = "<NOP>"

# Row 42: price option page - filled mandatory field
$ws.Range("A42").Value = "Price option page check for filled mandatory field"
$ws.Range("B42").Value = "<CHK>"
$ws.Range("F42").Value = "Price option page check for filled mandatory field"
$ws.Range("H42").Value = "<NOP>"

# --- Column F width (widened / best-fit to the new, longer content) --------
$ws.Columns("F").AutoFit()

# --- D1 no longer carries the stray "applyFill" style -----------------------
$ws.Range("D1").ClearFormats()

# --- Picture shifts down with the new rows ----------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = 631.8000787401575

# --- Selection / view state --------------------------------------------------
$ws.Activate()
$ws.Range("G28").Select()
